# Added periodic & upfront related scenarios
$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Update the value in B17 of the input sheet
# (was "RBI (India)", now the new scenario value)
$wsInput.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Make ProductLoanInput the active (selected) sheet with B17 as the active cell
$wsInput.Activate()
$wsInput.Range("B17").Select()

$wb.Save()
